$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row2 = $t.Rows.Item(2)
$cell = $row2.Cells.Item(1)
$newTableXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:tblPr><w:tblW w:w="14962" w:type="dxa"/><w:tblInd w:w="16" w:type="dxa"/><w:tblLayout w:type="fixed"/><w:tblCellMar><w:top w:w="28" w:type="dxa"/><w:left w:w="28" w:type="dxa"/><w:bottom w:w="28" w:type="dxa"/><w:right w:w="28" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2263"/><w:gridCol w:w="343"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="344"/><w:gridCol w:w="343"/><w:gridCol w:w="344"/><w:gridCol w:w="343"/><w:gridCol w:w="6854"/></w:tblGrid><w:tr w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w14:paraId="29572BA4" w14:textId="77777777" w:rsidTr="009D5FDC"><w:trPr><w:trHeight w:val="125"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2263" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="000000"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="639EE51A" w14:textId="77777777" w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Simplified Arabic Fixed"/><w:b/><w:noProof/><w:color w:val="FFFFFF"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r w:rsidRPr="00FF7E8F"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Simplified Arabic Fixed"/><w:b/><w:noProof/><w:color w:val="FFFFFF"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="000000"/><w:lang w:val="id-ID"/></w:rPr><w:t>NO PENDAFT</w:t></w:r><w:r w:rsidRPr="00FF7E8F"><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Simplified Arabic Fixed"/><w:b/><w:noProof/><w:color w:val="FFFFFF"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>ARAN</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="343" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="61C7FDFF" w14:textId="77777777" w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="2EB2FBE4" w14:textId="6D6FB835" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="2C5EEB11" w14:textId="76AAC319" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="58DF87D9" w14:textId="100B8833" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="76B2DDD7" w14:textId="2574D895" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="1C83CC44" w14:textId="0EA58D9F" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn5</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="60616963" w14:textId="0036B6CD" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn6</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="0CC41FE7" w14:textId="4E0F8620" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn7</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="41CEB43B" w14:textId="4499678C" w:rsidR="009D5FDC" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidR="00C64A4F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no_npsn8</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>_jalur</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>_jalur2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>_jalur3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>_jalur4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="5067FA86" w14:textId="3C69CC48" w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="343" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="4002CD2F" w14:textId="77777777" w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="344" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="46864A22" w14:textId="77777777" w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="343" w:type="dxa"/><w:tcBorders><w:top w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="28D6CDD3" w14:textId="77777777" w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/></w:rPr><w:t>no4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="6854" w:type="dxa"/><w:tcBorders><w:top w:val="nil"/><w:left w:val="dotted" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="nil"/><w:right w:val="dotted" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:p w14:paraId="53A2D877" w14:textId="77777777" w:rsidR="009D5FDC" w:rsidRPr="00FF7E8F" w:rsidRDefault="009D5FDC" w:rsidP="001D1727"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow" w:cs="Simplified Arabic Fixed"/><w:b/><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p></w:tc></w:tr></w:tbl>'
$cell.Range.InsertXML($newTableXml)
Write-Host "Done"
